$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '40.536.58'
$ws.Range('E2').Value = '  -2.50%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.369.38'
$ws.Range('E3').Value = '  -3.94%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.50'
$ws.Range('E5').Value = '  -2.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '86.86'
$ws.Range('E6').Value = '  -5.74%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.528'
$ws.Range('E7').Value = '  -4.17%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.495'
$ws.Range('E9').Value = '  -3.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0842'
$ws.Range('E10').Value = '  -2.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '30.60'
$ws.Range('E11').Value = '  -6.81%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.737.77'
$ws.Range('E13').Value = '  -3.82%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.53'
$ws.Range('E14').Value = '  -4.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.08'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.368.07'
$ws.Range('E16').Value = '  -4.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.758'
$ws.Range('E17').Value = '  -4.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '40.508.03'
$ws.Range('E18').Value = '  -2.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0911'
$ws.Range('E19').Value = '  -3.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.13'
$ws.Range('E20').Value = '  -4.53%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '68.42'
$ws.Range('E21').Value = '  -3.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.74'
$ws.Range('E22').Value = '  -4.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.16'
$ws.Range('E23').Value = '  -1.98%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.58'
$ws.Range('E24').Value = '  -6.15%  '
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.80'
$ws.Range('E26').Value = '  -7.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.76'
$ws.Range('E27').Value = '  -3.91%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.20'
$ws.Range('E28').Value = '  -2.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.30'
$ws.Range('E29').Value = '  -3.81%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.19'
$ws.Range('E30').Value = '  -5.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '154.90'
$ws.Range('E31').Value = '  -1.40%  '
$ws.Range('E32').Value = '  +0.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.21'
$ws.Range('E33').Value = '  -4.34%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0728'
$ws.Range('E34').Value = '  -4.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.41'
$ws.Range('E35').Value = '  -5.71%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.113'
$ws.Range('E36').Value = '  -2.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '16.11'
$ws.Range('E37').Value = '  -6.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.77'
$ws.Range('E38').Value = '  -4.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0998'
$ws.Range('E39').Value = '  -3.97%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.70'
$ws.Range('E40').Value = '  -7.79%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.81'
$ws.Range('E41').Value = '  -4.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.38'
$ws.Range('E42').Value = '  -6.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.966.36'
$ws.Range('E43').Value = '  -1.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0269'
$ws.Range('E44').Value = '  -4.97%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.84'
$ws.Range('E45').Value = '  -5.30%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.40'
$ws.Range('E46').Value = '  -0.28%  '
$ws.Range('E47').Value = '  -7.99%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.601.89'
$ws.Range('E48').Value = '  -3.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '93.02'
$ws.Range('E49').Value = '  -4.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.25'
$ws.Range('E50').Value = '  -3.89%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '50.31'
$ws.Range('E51').Value = '  -3.46%  '
